$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the two new columns (I0, IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the header text (PasteSpecial of formats only should not touch values, but set again to be safe)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0) and column J (IF), rows 2-74
$iValues = @(9,9,9,6,6,7,4,7,10,8,5,6,4,8,6,6,7,5,8,5,7,7,7,8,3,8,7,6,7,7,7,6,9,9,7,7,4,6,7,8,7,9,9,9,5,6,7,9,10,8,9,8,5,8,5,9,7,8,6,7,7,5,7,8,7,6,9,9,7,5,7,6,3)
$jValues = @(9,9,9,6,6,7,6,7,10,9,5,6,4,8,6,7,8,6,9,6,7,7,7,8,4,8,8,7,8,8,8,6,9,9,7,8,6,8,8,9,8,9,9,9,6,8,8,9,10,8,9,8,6,8,6,9,7,8,6,7,7,6,8,8,8,7,9,9,8,5,8,6,3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}

